# Sprint backlog rewrite: convert the flat Product Backlog into Scrum sprint backlogs,
# add a "Sprint" table column, and resize the table / column widths accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old row 10 content (its data shifted up into row 9 in the new layout) ---
$ws.Range("A10:D10").ClearContents()

# --- Extend the existing table with a new "Sprint" column ---
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null

# --- Populate header + data rows (A: Problem, B: ID, C: Estimation, D: Priority, E: Sprint, G: Sprint legend) ---
# Row 1
$ws.Range("A1").Value2 = "Problem"
$ws.Range("B1").Value2 = "ID"
$ws.Range("C1").Value2 = "Estimation (Time for task)"
$ws.Range("D1").Value2 = "Priority"
$ws.Range("E1").Value2 = "Sprint"

# Row 2
$ws.Range("A2").Value2 = "Group formation"
$ws.Range("B2").Value2 = 1
$ws.Range("C2").Value2 = "N/A"
$ws.Range("D2").Value2 = 1
$ws.Range("E2").Value2 = "1st"
$ws.Range("G2").Value2 = "1st Sprint - weeks 1-3"

# Row 3
$ws.Range("A3").Value2 = "Git Hub"
$ws.Range("B3").Value2 = 2
$ws.Range("C3").Value2 = "N/A"
$ws.Range("D3").Value2 = 2
$ws.Range("E3").Value2 = "1st"
$ws.Range("G3").Value2 = "2nd Sprint - weeks 4-6"

# Row 4
$ws.Range("A4").Value2 = "Code of Conduct"
$ws.Range("B4").Value2 = 3
$ws.Range("C4").Value2 = "Ten mintues"
$ws.Range("D4").Value2 = 8
$ws.Range("E4").Value2 = "1st"
$ws.Range("G4").Value2 = "3rd Sprint - weeks 7-9"

# Row 5
$ws.Range("A5").Value2 = "Travis CI"
$ws.Range("B5").Value2 = 4
$ws.Range("C5").Value2 = "2 weeks"
$ws.Range("D5").Value2 = 3
$ws.Range("E5").Value2 = "1st"
$ws.Range("G5").Value2 = "4th Sprint - weeks 10-12"

# Row 6
$ws.Range("A6").Value2 = "Docker implmentation"
$ws.Range("B6").Value2 = 5
$ws.Range("C6").Value2 = "1 day"
$ws.Range("D6").Value2 = 4
$ws.Range("E6").Value2 = "1st"

# Row 7
$ws.Range("A7").Value2 = "JAR"
$ws.Range("B7").Value2 = 6
$ws.Range("D7").Value2 = 5
$ws.Range("E7").Value2 = "1st"

# Row 8
$ws.Range("A8").Value2 = "Maven"
$ws.Range("B8").Value2 = 7
$ws.Range("C8").Value2 = "1 day"
$ws.Range("D8").Value2 = 6
$ws.Range("E8").Value2 = "1st"

# Row 9
$ws.Range("A9").Value2 = "First Release"
$ws.Range("B9").Value2 = 8
$ws.Range("C9").Value2 = "Three Weeks"
$ws.Range("D9").Value2 = 7
$ws.Range("E9").Value2 = "1st"

# Row 11
$ws.Range("A11").Value2 = "Integration with Zube.io"
$ws.Range("B11").Value2 = 9
$ws.Range("C11").Value2 = "1 day"
$ws.Range("D11").Value2 = 1
$ws.Range("E11").Value2 = "2nd"

# Row 12
$ws.Range("A12").Value2 = "UML diagram creation"
$ws.Range("B12").Value2 = 10
$ws.Range("C12").Value2 = "2 days"
$ws.Range("D12").Value2 = 2
$ws.Range("E12").Value2 = "2nd"

# Row 13
$ws.Range("A13").Value2 = "Introduction of sprint boards"
$ws.Range("B13").Value2 = 11
$ws.Range("C13").Value2 = "1 day"
$ws.Range("D13").Value2 = 3
$ws.Range("E13").Value2 = "2nd"

# Row 14
$ws.Range("A14").Value2 = "Preparation of first 10 raports"
$ws.Range("B14").Value2 = 12
$ws.Range("C14").Value2 = "2 weeks"
$ws.Range("D14").Value2 = 4
$ws.Range("E14").Value2 = "2nd"

# Row 16
$ws.Range("A16").Value2 = "Introduction of TDD"
$ws.Range("B16").Value2 = 13
$ws.Range("C16").Value2 = "1 day"
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = "3rd"

# Row 17
$ws.Range("A17").Value2 = "Unit tests"
$ws.Range("B17").Value2 = 14
$ws.Range("C17").Value2 = "1 week"
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = "3rd"

# Row 18
$ws.Range("A18").Value2 = "Integrating tests with Travis"
$ws.Range("B18").Value2 = 15
$ws.Range("C18").Value2 = "1 day"
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = "3rd"

# Row 19
$ws.Range("A19").Value2 = "Preparation of remaining 15 raports"
$ws.Range("B19").Value2 = 16
$ws.Range("C19").Value2 = "2 weeks"
$ws.Range("D19").Value2 = 2
$ws.Range("E19").Value2 = "3rd"

# Row 21
$ws.Range("A21").Value2 = "Bug reporting system"
$ws.Range("B21").Value2 = 17
$ws.Range("C21").Value2 = "1 day"
$ws.Range("D21").Value2 = 1
$ws.Range("E21").Value2 = "4th"

# Row 22
$ws.Range("A22").Value2 = "Deployment working"
$ws.Range("B22").Value2 = 19
$ws.Range("C22").Value2 = "3 weeks"
$ws.Range("D22").Value2 = 2
$ws.Range("E22").Value2 = "4th"

# --- Right-align the Sprint column data cells (incl. the blank spacer row 15) ---
$ws.Range("E2").HorizontalAlignment = -4152
$ws.Range("E3").HorizontalAlignment = -4152
$ws.Range("E4").HorizontalAlignment = -4152
$ws.Range("E5").HorizontalAlignment = -4152
$ws.Range("E6").HorizontalAlignment = -4152
$ws.Range("E7").HorizontalAlignment = -4152
$ws.Range("E8").HorizontalAlignment = -4152
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E12").HorizontalAlignment = -4152
$ws.Range("E13").HorizontalAlignment = -4152
$ws.Range("E14").HorizontalAlignment = -4152
$ws.Range("E15").HorizontalAlignment = -4152
$ws.Range("E16").HorizontalAlignment = -4152
$ws.Range("E17").HorizontalAlignment = -4152
$ws.Range("E18").HorizontalAlignment = -4152
$ws.Range("E19").HorizontalAlignment = -4152
$ws.Range("E21").HorizontalAlignment = -4152
$ws.Range("E22").HorizontalAlignment = -4152

# --- Resize the table to cover the full (pre-enlarged) range used by the author ---
$lo.Resize($ws.Range("A1:E41"))

# --- Column widths: A, B widened for new text; G added for the sprint legend ---
$ws.Columns.Item(1).ColumnWidth = 33.166666666666664   # -> stored width 34
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666   # -> stored width closest to 11.5703125
$ws.Columns.Item(7).ColumnWidth = 26.833333333333332   # -> stored width closest to 27.7109375

# --- Final selection, matching the saved cursor position ---
$ws.Range("D22").Select()